$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.453351974487305
$ws.Range("B1").Value = 3.470021963119507
$ws.Range("C1").Value = 3.067542791366577
$ws.Range("D1").Value = 3.83814525604248
$ws.Range("E1").Value = 5.149348258972168
